# Updates cryptos list prices/volume percentages (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.899.37'
$ws.Cells.Item(2, 5).Value = '  -0.28%  '

$ws.Cells.Item(3, 4).Value = '2.305.52'
$ws.Cells.Item(3, 5).Value = '  -0.78%  '

$ws.Cells.Item(4, 5).Value = '  +0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '108.87'
$ws.Cells.Item(5, 5).Value = '  +11.72%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '271.68'
$ws.Cells.Item(6, 5).Value = '  -0.14%  '

$ws.Cells.Item(7, 5).Value = '  -0.60%  '

$ws.Cells.Item(8, 5).Value = '  +0.16%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.617'
$ws.Cells.Item(9, 5).Value = '  -1.74%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '47.58'
$ws.Cells.Item(10, 5).Value = '  +4.52%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0940'
$ws.Cells.Item(11, 5).Value = '  -1.65%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '8.47'
$ws.Cells.Item(12, 5).Value = '  +5.56%  '

$ws.Cells.Item(13, 5).Value = '  +1.53%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '15.78'
$ws.Cells.Item(14, 5).Value = '  +1.67%  '

$ws.Cells.Item(15, 4).Value = '2.647.77'
$ws.Cells.Item(15, 5).Value = '  -0.28%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.862'
$ws.Cells.Item(16, 5).Value = '  -1.64%  '

$ws.Cells.Item(17, 4).Value = '2.300.82'
$ws.Cells.Item(17, 5).Value = '  -0.95%  '

$ws.Cells.Item(18, 4).Value = '43.803.89'
$ws.Cells.Item(18, 5).Value = '  -0.35%  '

$ws.Cells.Item(19, 5).Value = '  +1.66%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.31'
$ws.Cells.Item(20, 5).Value = '  -1.70%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '72.36'
$ws.Cells.Item(21, 5).Value = '  -1.68%  '

$ws.Cells.Item(22, 5).Value = '  +8.05%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '234.66'
$ws.Cells.Item(23, 5).Value = '  -2.37%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.98'
$ws.Cells.Item(24, 5).Value = '  +17.42%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.25'
$ws.Cells.Item(25, 5).Value = '  -1.70%  '

$ws.Cells.Item(26, 5).Value = '  +0.00%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.40'
$ws.Cells.Item(27, 5).Value = '  +0.06%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '41.19'
$ws.Cells.Item(28, 5).Value = '  +7.68%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '3.45'
$ws.Cells.Item(29, 5).Value = '  -1.43%  '

$ws.Cells.Item(30, 5).Value = '  -0.99%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '177.76'
$ws.Cells.Item(31, 5).Value = '  +1.25%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '21.92'
$ws.Cells.Item(32, 5).Value = '  -2.39%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0916'
$ws.Cells.Item(33, 5).Value = '  +0.26%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.62'
$ws.Cells.Item(34, 5).Value = '  +2.06%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.81'
$ws.Cells.Item(35, 5).Value = '  +7.71%  '

$ws.Cells.Item(36, 5).Value = '  -0.47%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.115'
$ws.Cells.Item(37, 5).Value = '  +4.10%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0357'
$ws.Cells.Item(38, 5).Value = '  -2.30%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.81'
$ws.Cells.Item(39, 5).Value = '  +12.52%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.238'
$ws.Cells.Item(40, 5).Value = '  -3.23%  '

$ws.Cells.Item(41, 5).Value = '  -1.92%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.39'
$ws.Cells.Item(42, 5).Value = '  -2.14%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '67.13'
$ws.Cells.Item(43, 5).Value = '  +6.71%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '12.21'
$ws.Cells.Item(44, 5).Value = '  -1.87%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '5.53'
$ws.Cells.Item(45, 5).Value = '  +3.56%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '8.81'
$ws.Cells.Item(46, 5).Value = '  -3.38%  '

$ws.Cells.Item(47, 5).Value = '  -1.92%  '

$ws.Cells.Item(48, 5).Value = '  +1.73%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '99.54'
$ws.Cells.Item(49, 5).Value = '  -0.96%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.439'
$ws.Cells.Item(50, 5).Value = '  +4.87%  '

$ws.Cells.Item(51, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(51, 4).Value = '2.537.64'
$ws.Cells.Item(51, 5).Value = '  -0.30%  '
